$d = $word.ActiveDocument

# Locate the paragraph that currently holds the single professor entry
# ("5111420 - Talita Martins Lacerda") under "Docente(s) Responsável(eis)".
$rng = $d.Content
$rng.Find.ClearFormatting()
$found = $rng.Find.Execute("5111420 - Talita Martins Lacerda", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find target paragraph text"
}

$para = $rng.Paragraphs(1)
$target = $para.Range

# Build the replacement paragraph content: six professor names, each its own
# run terminated with a manual line break except for the last one, matching
# the pattern already used elsewhere in the document (e.g. "Créditos-aula").
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
    '<pkg:xmlData>' + `
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
    '<w:body>' + `
    '<w:p>' + `
    '<w:pPr><w:pStyle w:val="ListBullet"/></w:pPr>' + `
    '<w:r><w:t>2143261 - André Luis Ferraz</w:t><w:br/></w:r>' + `
    '<w:r><w:t>5082401 - André Moreni Lopes</w:t><w:br/></w:r>' + `
    '<w:r><w:t>4873328 - Fernando Segato</w:t><w:br/></w:r>' + `
    '<w:r><w:t>6007846 - Júlio César dos Santos</w:t><w:br/></w:r>' + `
    '<w:r><w:t>5111420 - Talita Martins Lacerda</w:t><w:br/></w:r>' + `
    '<w:r><w:t>5817181 - Valdeir Arantes</w:t></w:r>' + `
    '</w:p>' + `
    '</w:body>' + `
    '</w:document>' + `
    '</pkg:xmlData>' + `
    '</pkg:part>' + `
    '</pkg:package>'

$target.InsertXML($xml)
